# Updated cryptos list on Thu Feb 23 06:57:44 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to be stored as literal text (not auto-converted to a number),
# while leaving the cell's style/format unaffected (reset back to Normal afterward).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "24.422.37"
$ws.Range("E2").Value = "  +1.60%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.666.69"
$ws.Range("E3").Value = "  +1.68%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -0.42%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "313.30"
$ws.Range("E5").Value = "  +1.91%  "

# Row 6 - USDC
Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  -0.14%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.3947"
$ws.Range("E7").Value = "  +1.26%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.3913"
$ws.Range("E8").Value = "  +1.73%  "

# Row 9 - OKB
Set-TextValue $ws.Range("D9") "52.11"
$ws.Range("E9").Value = "  +6.34%  "

# Row 10 - Polygon (price unchanged)
$ws.Range("E10").Value = "  +3.62%  "

# Row 11 - BinanceUSD (price unchanged)
$ws.Range("E11").Value = "  -0.43%  "

# Row 12 - Dogecoin
Set-TextValue $ws.Range("D12") "0.08576"
$ws.Range("E12").Value = "  +1.54%  "

# Row 13 - Solana
Set-TextValue $ws.Range("D13") "24.36"
$ws.Range("E13").Value = "  +2.04%  "

# Row 14 - Polkadot (price unchanged)
$ws.Range("E14").Value = "  +2.50%  "

# Row 15 - Chainlink
Set-TextValue $ws.Range("D15") "7.943"
$ws.Range("E15").Value = "  +6.37%  "

# Row 16 - ShibaInu (price unchanged)
$ws.Range("E16").Value = "  +5.02%  "

# Row 17 - WrappedEther
Set-TextValue $ws.Range("D17") "1.661.03"
$ws.Range("E17").Value = "  +1.10%  "

# Row 18 - Litecoin
Set-TextValue $ws.Range("D18") "95.05"
$ws.Range("E18").Value = "  +0.71%  "

# Row 19 - TRON
Set-TextValue $ws.Range("D19") "0.06994"
$ws.Range("E19").Value = "  +0.85%  "

# Row 20 - Avalanche
Set-TextValue $ws.Range("D20") "20.55"
$ws.Range("E20").Value = "  -0.96%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "6.987"
$ws.Range("E21").Value = "  +1.20%  "

# Row 22 - Dai (price unchanged)
$ws.Range("E22").Value = "  -0.27%  "

# Row 23 - Cosmos (price unchanged)
$ws.Range("E23").Value = "  +1.01%  "

# Row 24 - WrappedBTC
Set-TextValue $ws.Range("D24") "24.404.90"
$ws.Range("E24").Value = "  +1.49%  "

# Row 25 - Toncoin
Set-TextValue $ws.Range("D25") "2.427"
$ws.Range("E25").Value = "  +3.81%  "

# Row 26 - LidoDAOToken
Set-TextValue $ws.Range("D26") "3.043"
$ws.Range("E26").Value = "  +13.97%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "22.51"
$ws.Range("E27").Value = "  +0.50%  "

# Row 28 - Monero
Set-TextValue $ws.Range("D28") "157.22"
$ws.Range("E28").Value = "  -0.34%  "

# Rows 29 & 30 swap places: BitcoinCash <-> HuobiToken (with updated price/volume)
$ws.Range("B29").Value = "HuobiToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D29") "5.464"
$ws.Range("E29").Value = "  +3.29%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D30") "142.73"
$ws.Range("E30").Value = "  +0.57%  "

# Row 31 - Filecoin
Set-TextValue $ws.Range("D31") "8.047"
$ws.Range("E31").Value = "  -8.73%  "

# Row 32 - WEMIXTOKEN (price unchanged)
$ws.Range("E32").Value = "  +3.21%  "

# Row 33 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D33") "1.842.65"
$ws.Range("E33").Value = "  +0.87%  "

# Row 34 - ImmutableX
Set-TextValue $ws.Range("D34") "1.058"
$ws.Range("E34").Value = "  +9.71%  "

# Row 35 - Hedera
Set-TextValue $ws.Range("D35") "0.08236"
$ws.Range("E35").Value = "  +2.80%  "

# Row 36 - VeChain
Set-TextValue $ws.Range("D36") "0.03031"
$ws.Range("E36").Value = "  +4.04%  "

# Row 37 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D37") "6.897"
$ws.Range("E37").Value = "  -3.84%  "

# Row 38 - FraxShare
Set-TextValue $ws.Range("D38") "11.12"
$ws.Range("E38").Value = "  +11.43%  "

# Row 39 - Algorand
Set-TextValue $ws.Range("D39") "0.2763"
$ws.Range("E39").Value = "  +2.38%  "

# Row 40 - Stellar (price unchanged)
$ws.Range("E40").Value = "  +0.12%  "

# Row 41 - TheSandbox
Set-TextValue $ws.Range("D41") "0.7712"
$ws.Range("E41").Value = "  +1.51%  "

# Row 42 - Aptos
Set-TextValue $ws.Range("D42") "13.80"
$ws.Range("E42").Value = "  +5.89%  "

# Row 43 - TrustWalletToken
Set-TextValue $ws.Range("D43") "1.449"
$ws.Range("E43").Value = "  -0.82%  "

# Row 44 - EnergySwap
Set-TextValue $ws.Range("D44") "16.51"
$ws.Range("E44").Value = "  +3.51%  "

# Row 45 - Decentraland
Set-TextValue $ws.Range("D45") "0.7102"
$ws.Range("E45").Value = "  +3.38%  "

# Row 46 - NEARProtocol
Set-TextValue $ws.Range("D46") "2.534"
$ws.Range("E46").Value = "  +2.37%  "

# Row 47 - PancakeSwap
Set-TextValue $ws.Range("D47") "4.132"
$ws.Range("E47").Value = "  +1.19%  "

# Row 48 - Frax (volume unchanged)
Set-TextValue $ws.Range("D48") "1.000"

# Row 49 - Cronos
Set-TextValue $ws.Range("D49") "0.08428"
$ws.Range("E49").Value = "  +0.62%  "

# Row 50 - Quant
Set-TextValue $ws.Range("D50") "136.76"
$ws.Range("E50").Value = "  +2.17%  "

# Row 51 - Flow
Set-TextValue $ws.Range("D51") "1.269"
$ws.Range("E51").Value = "  +1.03%  "
